# Update relay settings (F1_relay_banshee) with the latest Banshee values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

# Row 3 (Relay #2) had stale slip settings that didn't match the rest of the
# table - bring it in line with the other relays.
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 0.5

# CT Primary (col F) for every relay moves from the placeholder 0.5 to the
# real ratio value of 13.
$ws.Range("F2:F15").Value = 13

# 27P Trip Pickup [pu] drops from 0.8 to 0.7 for every relay.
$ws.Range("N2:N15").Value = 0.7

# 50P IOC Trip Delay [Sec] (col Q) moves from 1 to 0 for every relay.
$ws.Range("Q2:Q15").Value = 0

# 51P Curve Type (col S) moves from 2 to 1 for every relay.
$ws.Range("S2:S15").Value = 1

# 51P TOC Time Dial (col T) moves from 3 to 1 for every relay.
$ws.Range("T2:T15").Value = 1

# 27P Trip Pickup [pu] (col U) moves from 0.9 to 0.5 for every relay.
$ws.Range("U2:U15").Value = 0.5

# 59P Trip Pickup [pu] (col V) moves from 1.1 to 1.2 for every relay.
$ws.Range("V2:V15").Value = 1.2

# 51P TOC Trip Pickup recalculated for relay #2 (row 3) after the slip
# settings fix above.
$ws.Range("K3").Value = 1.4701175847416352

# The "51P TOC Trip Pickup" duplicate column R is kept in sync with column K
# for every relay - fix the rows where it had drifted out of sync.
$ws.Range("R3").Value = 1.4701175847416352
$ws.Range("R10").Value = 1.4128571428571426
$ws.Range("R12").Value = 1.4128571428571426
$ws.Range("R14").Value = 1.4128571428571426
